$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-27 Friday" "2025-06-28 Saturday"

Replace-Text "677×6=" "631×5="
Replace-Text "466×6=" "989×4="
Replace-Text "895×4=" "300×9="
Replace-Text "630×8=" "115×4="
Replace-Text "492×6=" "678×6="
Replace-Text "104×9=" "269×6="
Replace-Text "988×4=" "574×8="
Replace-Text "698×7=" "932×7="
Replace-Text "150×4=" "944×7="
Replace-Text "858×2=" "400×9="
Replace-Text "186×8=" "946×7="
Replace-Text "164×8=" "531×6="
Replace-Text "816×3=" "461×2="
Replace-Text "923×7=" "536×9="
Replace-Text "379×7=" "718×5="
Replace-Text "443×5=" "929×2="
Replace-Text "532×2=" "388×8="
Replace-Text "828×7=" "713×2="
Replace-Text "874×8=" "522×4="
Replace-Text "338×2=" "890×9="
Replace-Text "664×4=" "362×2="
Replace-Text "548×5=" "366×4="
Replace-Text "678×4=" "637×5="
Replace-Text "887×4=" "498×2="
Replace-Text "915×9=" "800×2="
